$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1176.6842
$ws.Range("J17").Value = 1176.6842
$ws.Range("L17").Value = 3530.0526
$ws.Range("N17").Value = -3866.0526
$ws.Range("H19").Value = 1718
$ws.Range("I19").Value = 2500
$ws.Range("K19").Value = 2500
$ws.Range("M19").Value = -2325
$ws.Range("H58").Value = 654646.5
$ws.Range("I58").Value = 1089648.8
$ws.Range("J58").Value = 2143.0833
$ws.Range("K58").Value = 3268946.4
$ws.Range("L58").Value = 6429.249899999999
$ws.Range("M58").Value = -3268796.4
$ws.Range("N58").Value = -6729.249899999999
$ws.Range("H92").Value = 698.13336
$ws.Range("I92").Value = 737.3570999999999
$ws.Range("J92").Value = 149
$ws.Range("K92").Value = 737.3570999999999
$ws.Range("L92").Value = 149
$ws.Range("M92").Value = 510.6429000000001
$ws.Range("N92").Value = -2645
$ws.Range("H125").Value = 3413.2
$ws.Range("I125").Value = 5532
$ws.Range("J125").Value = 3087.2307
$ws.Range("K125").Value = 49788
$ws.Range("L125").Value = 27785.0763
$ws.Range("M125").Value = -47328
$ws.Range("N125").Value = -32705.0763
$ws.Range("H128").Value = 47918.285
$ws.Range("I128").Value = 43000
$ws.Range("J128").Value = 48738
$ws.Range("K128").Value = 43000
$ws.Range("L128").Value = 48738
$ws.Range("M128").Value = -38020
$ws.Range("N128").Value = -58698
$ws.Range("H137").Value = 1169.2858
$ws.Range("I137").Value = 1134.9524
$ws.Range("K137").Value = 3404.857199999999
$ws.Range("M137").Value = -854.8571999999995
$ws.Range("H141").Value = 2670.3333
$ws.Range("I141").Value = 2060
$ws.Range("K141").Value = 6180
$ws.Range("M141").Value = -1000
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 860.5599999999999
$ws.Range("I74").Value = 773.8946999999999
$ws.Range("J74").Value = 1135
$ws.Range("K74").Value = 773.8946999999999
$ws.Range("L74").Value = 1135
$ws.Range("M74").Value = 100.1053000000001
$ws.Range("N74").Value = -2883
$ws.Range("H77").Value = 860.5599999999999
$ws.Range("I77").Value = 773.8946999999999
$ws.Range("J77").Value = 1135
$ws.Range("K77").Value = 3869.4735
$ws.Range("L77").Value = 5675
$ws.Range("M77").Value = 498.5265000000004
$ws.Range("N77").Value = -14411
$ws.Range("H110").Value = 66807096
$ws.Range("I110").Value = 71578990
$ws.Range("J110").Value = 500
$ws.Range("K110").Value = 71578990
$ws.Range("L110").Value = 500
$ws.Range("M110").Value = -71576945
$ws.Range("N110").Value = -4590
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2000
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1549
$ws.Range("N94").ClearContents()
$ws.Range("H99").Value = 2152.5
$ws.Range("I99").Value = 2515
$ws.Range("J99").Value = 2112.2222
$ws.Range("K99").Value = 2515
$ws.Range("L99").Value = 2112.2222
$ws.Range("M99").Value = -1017
$ws.Range("N99").Value = -5108.2222
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2613.4
$ws.Range("I58").Value = 2423.9167
$ws.Range("J58").Value = 3371.3333
$ws.Range("K58").Value = 2423.9167
$ws.Range("L58").Value = 3371.3333
$ws.Range("M58").Value = -2220.9167
$ws.Range("N58").Value = -3777.3333
$ws.Range("H99").Value = 2577.4
$ws.Range("I99").Value = 3980
$ws.Range("J99").Value = 2421.5557
$ws.Range("K99").Value = 3980
$ws.Range("L99").Value = 2421.5557
$ws.Range("M99").Value = -2482
$ws.Range("N99").Value = -5417.5557
$ws.Range("H105").Value = 1446.5
$ws.Range("I105").Value = 1233
$ws.Range("J105").Value = 1660
$ws.Range("K105").Value = 1233
$ws.Range("L105").Value = 1660
$ws.Range("M105").Value = 514
$ws.Range("N105").Value = -5154
$ws.Range("H109").Value = 38000
$ws.Range("J109").Value = 38000
$ws.Range("L109").Value = 38000
$ws.Range("N109").Value = -40080
$ws.Range("H126").Value = 2577.4
$ws.Range("I126").Value = 3980
$ws.Range("J126").Value = 2421.5557
$ws.Range("K126").Value = 11940
$ws.Range("L126").Value = 7264.6671
$ws.Range("M126").Value = -9470
$ws.Range("N126").Value = -12204.6671
$ws.Range("H136").Value = 2613.4
$ws.Range("I136").Value = 2423.9167
$ws.Range("J136").Value = 3371.3333
$ws.Range("K136").Value = 7271.750100000001
$ws.Range("L136").Value = 10113.9999
$ws.Range("M136").Value = -4721.750100000001
$ws.Range("N136").Value = -15213.9999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 847.3333
$ws.Range("J2").Value = 1441.5714
$ws.Range("L2").Value = 8649.428400000001
$ws.Range("N2").Value = -8875.428400000001
$ws.Range("H3").Value = 3500
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 6000
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = -2888
$ws.Range("N3").Value = -18224
$ws.Range("H10").Value = 330
$ws.Range("I10").Value = 132.5
$ws.Range("K10").Value = 397.5
$ws.Range("M10").Value = -258.5
$ws.Range("H15").Value = 96.56521600000001
$ws.Range("I15").Value = 31.333334
$ws.Range("J15").Value = 218.875
$ws.Range("K15").Value = 94.00000199999999
$ws.Range("L15").Value = 656.625
$ws.Range("M15").Value = 45.99999800000001
$ws.Range("N15").Value = -936.625
$ws.Range("H21").Value = 100
$ws.Range("I21").Value = 100
$ws.Range("K21").Value = 300
$ws.Range("M21").Value = -127
$ws.Range("H26").Value = 1447
$ws.Range("I26").Value = 292
$ws.Range("J26").Value = 2602
$ws.Range("K26").Value = 876
$ws.Range("L26").Value = 7806
$ws.Range("M26").Value = -588
$ws.Range("N26").Value = -8382
$ws.Range("H32").Value = 999
$ws.Range("I32").Value = 999
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2997
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2714
$ws.Range("N32").ClearContents()
$ws.Range("H33").Value = 1898.1818
$ws.Range("I33").Value = 40
$ws.Range("J33").Value = 2084
$ws.Range("K33").Value = 240
$ws.Range("L33").Value = 12504
$ws.Range("M33").Value = 43
$ws.Range("N33").Value = -13070
$ws.Range("H34").Value = 1150.0769
$ws.Range("I34").Value = 470.8
$ws.Range("J34").Value = 1574.625
$ws.Range("K34").Value = 1412.4
$ws.Range("L34").Value = 4723.875
$ws.Range("M34").Value = -1328.4
$ws.Range("N34").Value = -4891.875
$ws.Range("H38").Value = 94.71429000000001
$ws.Range("J38").Value = 114.4
$ws.Range("L38").Value = 343.2
$ws.Range("N38").Value = -1037.2
$ws.Range("H40").Value = 746.8333
$ws.Range("I40").Value = 495.25
$ws.Range("J40").Value = 1250
$ws.Range("K40").Value = 1981
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -1912
$ws.Range("N40").Value = -5138
$ws.Range("H75").Value = 2577.1843
$ws.Range("I75").Value = 818.8333
$ws.Range("J75").Value = 2906.875
$ws.Range("K75").Value = 2456.4999
$ws.Range("L75").Value = 8720.625
$ws.Range("M75").Value = -1458.4999
$ws.Range("N75").Value = -10716.625
$ws.Range("H78").Value = 2577.1843
$ws.Range("I78").Value = 818.8333
$ws.Range("J78").Value = 2906.875
$ws.Range("K78").Value = 7369.4997
$ws.Range("L78").Value = 26161.875
$ws.Range("M78").Value = -2377.4997
$ws.Range("N78").Value = -36145.875
$ws.Range("H98").Value = 666832.3
$ws.Range("I98").Value = 500
$ws.Range("J98").Value = 999998.5
$ws.Range("K98").Value = 1500
$ws.Range("L98").Value = 2999995.5
$ws.Range("M98").Value = -2
$ws.Range("N98").Value = -3002991.5
$ws.Range("H105").Value = 7419.7
$ws.Range("J105").Value = 7688.5557
$ws.Range("L105").Value = 23065.6671
$ws.Range("N105").Value = -28307.6671
$ws.Range("H132").Value = 1803.6923
$ws.Range("I132").Value = 996.06665
$ws.Range("J132").Value = 2905
$ws.Range("K132").Value = 8964.599850000001
$ws.Range("L132").Value = 26145
$ws.Range("M132").Value = -6434.599850000001
$ws.Range("N132").Value = -31205
$ws.Range("H134").Value = 4234.227
$ws.Range("I134").Value = 3592.7144
$ws.Range("J134").Value = 4533.6
$ws.Range("K134").Value = 10778.1432
$ws.Range("L134").Value = 13600.8
$ws.Range("M134").Value = -5708.143199999999
$ws.Range("N134").Value = -23740.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 43000
$ws.Range("J105").Value = 43000
$ws.Range("L105").Value = 43000
$ws.Range("N105").Value = -49988
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 533297.4
$ws.Range("I46").Value = 330.9091
$ws.Range("J46").Value = 1266126.2
$ws.Range("K46").Value = 330.9091
$ws.Range("L46").Value = 1266126.2
$ws.Range("M46").Value = -142.9091
$ws.Range("N46").Value = -1266502.2
$ws.Range("H93").Value = 2163.0417
$ws.Range("I93").Value = 2229
$ws.Range("K93").Value = 2229
$ws.Range("M93").Value = -981
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 71429630
$ws.Range("I96").Value = 125001230
$ws.Range("J96").Value = 824
$ws.Range("K96").Value = 125001230
$ws.Range("L96").Value = 824
$ws.Range("M96").Value = -124999857
$ws.Range("N96").Value = -3570
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H122").Value = 2682.6843
$ws.Range("J122").Value = 3696.3333
$ws.Range("L122").Value = 11088.9999
$ws.Range("N122").Value = -15988.9999

Write-Output "Applied 252 changes"